# Applies the commit "Added URLs and views for all non-JSON pages":
#  1. Refresh the cached "Date Placeholder" text (01/08/2021 -> 10/08/2021)
#     on the slide master and every slide layout.
#  2. Update three URL-path fragments on slide 11's "Content Placeholder 2"
#     shape, turning trailing-slash-less paths into trailing-slash paths
#     (and collapsing "}/overview/" down to "}/"), without touching the
#     surrounding runs.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

function Replace-InTextRange($textRange, $search, $replace) {
    $full = $textRange.Text
    $idx = $full.IndexOf($search)
    if ($idx -ge 0) {
        $chars = $textRange.Characters($idx + 1, $search.Length)
        $chars.Text = $replace
    }
}

# --- 1. Date placeholder text on master + all layouts -----------------
$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes "10/08/2021"

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes "10/08/2021"
}

# --- 2. URL path fragments on slide 11 ---------------------------------
$slide = $p.Slides.Item(11)
$contentShape = $slide.Shapes.Item("Content Placeholder 2")
$tr = $contentShape.TextFrame.TextRange

Replace-InTextRange $tr "}/overview/" "}/"
Replace-InTextRange $tr "}/enter" "}/enter/"
Replace-InTextRange $tr "}/predict" "}/predict/"
